$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-unused bold style (s="2") from column M rows 151-172
$ws.Range("M151:M172").ClearFormats()

# Append new event rows 173-181
# Row 173
$ws.Range("A173").Value = '00:08:39.820'
$ws.Range("B173").Value = 'C.D. Retiro Sur'
$ws.Range("C173").Value = 'Pass'
$ws.Range("D173").Value = 8
$ws.Range("E173").Value = 68
$ws.Range("F173").Value = 25
$ws.Range("G173").Value = 84
$ws.Range("H173").Value = 32
$ws.Range("J173").Value = 'Incomplete'
$ws.Range("K173").Value = 'Ground Pass'
$ws.Range("L173").Value = 'Adrian Lopez'
$ws.Range("M173").Value = 'C.D. Retiro Sur'

# Row 174
$ws.Range("A174").Value = '00:08:40.820'
$ws.Range("B174").Value = 'Escuela Dep. Moratalaz ''D'''
$ws.Range("C174").Value = 'Interception'
$ws.Range("D174").Value = 8
$ws.Range("E174").Value = 84
$ws.Range("F174").Value = 32
$ws.Range("L174").Value = 'Fabio Patus'
$ws.Range("M174").Value = 'C.D. Retiro Sur'

# Row 175
$ws.Range("A175").Value = '00:08:41.820'
$ws.Range("B175").Value = 'Escuela Dep. Moratalaz ''D'''
$ws.Range("C175").Value = 'Ball Receipt'
$ws.Range("D175").Value = 8
$ws.Range("E175").Value = 84
$ws.Range("F175").Value = 42
$ws.Range("L175").Value = 'Diego Refoyo'
$ws.Range("M175").Value = 'Escuela Dep. Moratalaz ''D'''

# Row 176
$ws.Range("A176").Value = '00:08:42.820'
$ws.Range("B176").Value = 'Escuela Dep. Moratalaz ''D'''
$ws.Range("C176").Value = 'Pass'
$ws.Range("D176").Value = 8
$ws.Range("E176").Value = 84
$ws.Range("F176").Value = 42
$ws.Range("G176").Value = 74
$ws.Range("H176").Value = 61
$ws.Range("I176").Value = 'Adrian Pombo'
$ws.Range("J176").Value = 'Complete'
$ws.Range("K176").Value = 'Ground Pass'
$ws.Range("L176").Value = 'Diego Refoyo'
$ws.Range("M176").Value = 'Escuela Dep. Moratalaz ''D'''

# Row 177
$ws.Range("A177").Value = '00:08:47.820'
$ws.Range("B177").Value = 'Escuela Dep. Moratalaz ''D'''
$ws.Range("C177").Value = 'Pass'
$ws.Range("D177").Value = 8
$ws.Range("E177").Value = 62
$ws.Range("F177").Value = 70
$ws.Range("G177").Value = 76
$ws.Range("H177").Value = 60
$ws.Range("I177").Value = 'Ángel Jesús'
$ws.Range("J177").Value = 'Complete'
$ws.Range("K177").Value = 'High Pass'
$ws.Range("L177").Value = 'Adrian Pombo'
$ws.Range("M177").Value = 'Escuela Dep. Moratalaz ''D'''

# Row 178
$ws.Range("A178").Value = '00:08:51.820'
$ws.Range("B178").Value = 'Escuela Dep. Moratalaz ''D'''
$ws.Range("C178").Value = 'Pass'
$ws.Range("D178").Value = 8
$ws.Range("E178").Value = 76
$ws.Range("F178").Value = 60
$ws.Range("G178").Value = 95
$ws.Range("H178").Value = 30
$ws.Range("I178").Value = 'Pablo Escribano'
$ws.Range("J178").Value = 'Complete'
$ws.Range("K178").Value = 'Ground Pass'
$ws.Range("L178").Value = 'Ángel Jesús'
$ws.Range("M178").Value = 'Escuela Dep. Moratalaz ''D'''

# Row 179
$ws.Range("A179").Value = '00:08:56.820'
$ws.Range("B179").Value = 'Escuela Dep. Moratalaz ''D'''
$ws.Range("C179").Value = 'Pass'
$ws.Range("D179").Value = 8
$ws.Range("E179").Value = 81
$ws.Range("F179").Value = 18
$ws.Range("G179").Value = 65
$ws.Range("H179").Value = 3
$ws.Range("I179").Value = 'Carlos Enrique'
$ws.Range("J179").Value = 'Complete'
$ws.Range("K179").Value = 'Ground Pass'
$ws.Range("L179").Value = 'Pablo Escribano'
$ws.Range("M179").Value = 'Escuela Dep. Moratalaz ''D'''

# Row 180
$ws.Range("A180").Value = '00:09:00.820'
$ws.Range("B180").Value = 'Escuela Dep. Moratalaz ''D'''
$ws.Range("C180").Value = 'Pass'
$ws.Range("D180").Value = 9
$ws.Range("E180").Value = 65
$ws.Range("F180").Value = 3
$ws.Range("G180").Value = 53
$ws.Range("H180").Value = 9
$ws.Range("I180").Value = 'Alejandro Charro'
$ws.Range("J180").Value = 'Complete'
$ws.Range("K180").Value = 'Ground Pass'
$ws.Range("L180").Value = 'Carlos Enrique'
$ws.Range("M180").Value = 'Escuela Dep. Moratalaz ''D'''

# Row 181
$ws.Range("A181").Value = '00:09:02.820'
$ws.Range("B181").Value = 'Escuela Dep. Moratalaz ''D'''
$ws.Range("C181").Value = 'Pass'
$ws.Range("D181").Value = 9
$ws.Range("E181").Value = 52
$ws.Range("F181").Value = 8
$ws.Range("G181").Value = 50
$ws.Range("H181").Value = 2
$ws.Range("I181").Value = 'Santiago Sanchez'
$ws.Range("J181").Value = 'Complete'
$ws.Range("K181").Value = 'Ground Pass'
$ws.Range("L181").Value = 'Alejandro Charro'
$ws.Range("M181").Value = 'Escuela Dep. Moratalaz ''D'''


# Update the active selection / scroll position to match the saved view
$ws.Range("L20").Select()
